$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 125
$ws.Range("H125").Value = 2078
$ws.Range("J125").Value = 795
$ws.Range("L125").Value = 7155
$ws.Range("N125").Value = -12075
# Row 132
$ws.Range("H132").Value = 2527.139
$ws.Range("I132").Value = 2008.742
$ws.Range("K132").Value = 6026.226
$ws.Range("M132").Value = -3496.226
# Row 137
$ws.Range("H137").Value = 15089.4
$ws.Range("J137").Value = 18111.75
$ws.Range("L137").Value = 54335.25
$ws.Range("N137").Value = -59435.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1780.6666
$ws.Range("I2").Value = 1780.6666
$ws.Range("K2").Value = 1780.6666
$ws.Range("M2").Value = -1667.6666
# Row 26
$ws.Range("H26").Value = 3916.6667
$ws.Range("I26").Value = 3916.6667
$ws.Range("K26").Value = 3916.6667
$ws.Range("M26").Value = -3586.6667
# Row 31
$ws.Range("H31").Value = 58463
$ws.Range("I31").Value = 4999.6665
$ws.Range("J31").Value = 111926.336
$ws.Range("K31").Value = 4999.6665
$ws.Range("L31").Value = 111926.336
$ws.Range("M31").Value = -4705.6665
$ws.Range("N31").Value = -112514.336
# Row 45
$ws.Range("H45").Value = 3357.4
$ws.Range("J45").Value = 3874.75
$ws.Range("L45").Value = 3874.75
$ws.Range("N45").Value = -4628.75
# Row 61
$ws.Range("H61").Value = 20046510
$ws.Range("I61").Value = 27782960
$ws.Range("J61").Value = 152779.72
$ws.Range("K61").Value = 27782960
$ws.Range("L61").Value = 152779.72
$ws.Range("M61").Value = -27782748
$ws.Range("N61").Value = -153203.72
# Row 110
$ws.Range("H110").Value = 1974.2
$ws.Range("I110").Value = 1888
$ws.Range("K110").Value = 1888
$ws.Range("M110").Value = 157
# Row 116
$ws.Range("H116").Value = 1780.6666
$ws.Range("I116").Value = 1780.6666
$ws.Range("K116").Value = 1780.6666
$ws.Range("M116").Value = 513.3334
# Row 136
$ws.Range("H136").Value = 20046510
$ws.Range("I136").Value = 27782960
$ws.Range("J136").Value = 152779.72
$ws.Range("K136").Value = 83348880
$ws.Range("L136").Value = 458339.16
$ws.Range("M136").Value = -83346330
$ws.Range("N136").Value = -463439.16

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1780.6666
$ws.Range("I3").Value = 1780.6666
$ws.Range("K3").Value = 1780.6666
$ws.Range("M3").Value = -1666.6666
# Row 28
$ws.Range("H28").Value = 46955
$ws.Range("J28").Value = 46955
$ws.Range("L28").Value = 46955
$ws.Range("N28").Value = -47543
# Row 96
$ws.Range("H96").Value = 30463.555
$ws.Range("J96").Value = 65924.664
$ws.Range("L96").Value = 65924.664
$ws.Range("N96").Value = -71416.664
# Row 107
$ws.Range("H107").Value = 1471.5714
$ws.Range("I107").Value = 1236.6364
$ws.Range("J107").Value = 2333
$ws.Range("K107").Value = 1236.6364
$ws.Range("L107").Value = 2333
$ws.Range("M107").Value = 683.3635999999999
$ws.Range("N107").Value = -6173
# Row 119
$ws.Range("H119").Value = 63687.25
$ws.Range("J119").Value = 63687.25
$ws.Range("L119").Value = 63687.25
$ws.Range("N119").Value = -73363.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 2087.1177
$ws.Range("I7").Value = 137.13333
$ws.Range("J7").Value = 16712
$ws.Range("K7").Value = 137.13333
$ws.Range("L7").Value = 16712
$ws.Range("M7").Value = -24.13333
$ws.Range("N7").Value = -16938
# Row 33
$ws.Range("H33").Value = 4577
$ws.Range("I33").Value = 4577
$ws.Range("K33").Value = 4577
$ws.Range("M33").Value = -4198
# Row 58
$ws.Range("H58").Value = 1213.1111
$ws.Range("J58").Value = 2000.6666
$ws.Range("L58").Value = 2000.6666
$ws.Range("N58").Value = -2406.6666
# Row 80
$ws.Range("H80").Value = 69989.5
$ws.Range("J80").Value = 69989.5
$ws.Range("L80").Value = 69989.5
$ws.Range("N80").Value = -72235.5
# Row 83
$ws.Range("H83").Value = 69989.5
$ws.Range("J83").Value = 69989.5
$ws.Range("L83").Value = 209968.5
$ws.Range("N83").Value = -221200.5
# Row 86
$ws.Range("H86").Value = 749.5
$ws.Range("I86").Value = 500
$ws.Range("J86").Value = 999
$ws.Range("K86").Value = 500
$ws.Range("L86").Value = 999
$ws.Range("M86").Value = 623
$ws.Range("N86").Value = -3245
# Row 89
$ws.Range("H89").Value = 749.5
$ws.Range("I89").Value = 500
$ws.Range("J89").Value = 999
$ws.Range("K89").Value = 2500
$ws.Range("L89").Value = 4995
$ws.Range("M89").Value = 3116
$ws.Range("N89").Value = -16227
# Row 99
$ws.Range("H99").Value = 2592.2727
$ws.Range("I99").Value = 2185.9167
$ws.Range("J99").Value = 3079.9
$ws.Range("K99").Value = 2185.9167
$ws.Range("L99").Value = 3079.9
$ws.Range("M99").Value = -687.9167000000002
$ws.Range("N99").Value = -6075.9
# Row 103
$ws.Range("H103").Value = 34967.25
$ws.Range("J103").Value = 44935
$ws.Range("L103").Value = 44935
$ws.Range("N103").Value = -47279
# Row 105
$ws.Range("H105").Value = 1848.375
$ws.Range("I105").Value = 1197.125
$ws.Range("K105").Value = 1197.125
$ws.Range("M105").Value = 549.875
# Row 107
$ws.Range("H107").Value = 1443.8572
$ws.Range("I107").Value = 1456.4
$ws.Range("J107").Value = 1412.5
$ws.Range("K107").Value = 1456.4
$ws.Range("L107").Value = 1412.5
$ws.Range("M107").Value = 463.5999999999999
$ws.Range("N107").Value = -5252.5
# Row 126
$ws.Range("H126").Value = 2592.2727
$ws.Range("I126").Value = 2185.9167
$ws.Range("J126").Value = 3079.9
$ws.Range("K126").Value = 6557.750100000001
$ws.Range("L126").Value = 9239.700000000001
$ws.Range("M126").Value = -4087.750100000001
$ws.Range("N126").Value = -14179.7
# Row 132
$ws.Range("H132").Value = 2216.76
$ws.Range("I132").Value = 2101.261
$ws.Range("J132").Value = 3545
$ws.Range("K132").Value = 6303.782999999999
$ws.Range("L132").Value = 10635
$ws.Range("M132").Value = -3773.782999999999
$ws.Range("N132").Value = -15695
# Row 134
$ws.Range("H134").Value = 481122.56
$ws.Range("I134").Value = 771368.1
$ws.Range("J134").Value = 9473.5
$ws.Range("K134").Value = 2314104.3
$ws.Range("L134").Value = 28420.5
$ws.Range("M134").Value = -2311569.3
$ws.Range("N134").Value = -33490.5
# Row 136
$ws.Range("H136").Value = 1213.1111
$ws.Range("J136").Value = 2000.6666
$ws.Range("L136").Value = 6001.9998
$ws.Range("N136").Value = -11101.9998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 3832.7778
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 4061.875
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 12185.625
$ws.Range("M80").Value = -5064
$ws.Range("N80").Value = -14057.625
# Row 83
$ws.Range("H83").Value = 3832.7778
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 4061.875
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 36556.875
$ws.Range("M83").Value = -13320
$ws.Range("N83").Value = -45916.875
# Row 122
$ws.Range("H122").Value = 1312.85
$ws.Range("I122").Value = 699
$ws.Range("J122").Value = 1517.4667
$ws.Range("K122").Value = 6291
$ws.Range("L122").Value = 13657.2003
$ws.Range("M122").Value = -3841
$ws.Range("N122").Value = -18557.2003
# Row 131
$ws.Range("H131").Value = 5649.647
$ws.Range("I131").Value = 6983.222
$ws.Range("J131").Value = 4149.375
$ws.Range("K131").Value = 20949.666
$ws.Range("L131").Value = 12448.125
$ws.Range("M131").Value = -15909.666
$ws.Range("N131").Value = -22528.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 83666.336
$ws.Range("J15").Value = 83666.336
$ws.Range("L15").Value = 83666.336
$ws.Range("N15").Value = -84242.336
# Row 81
$ws.Range("H81").Value = 83666.336
$ws.Range("J81").Value = 83666.336
$ws.Range("L81").Value = 83666.336
$ws.Range("N81").Value = -85662.336
# Row 84
$ws.Range("H84").Value = 83666.336
$ws.Range("J84").Value = 83666.336
$ws.Range("L84").Value = 250999.008
$ws.Range("N84").Value = -260983.008

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3466.2856
$ws.Range("I40").Value = 2775.2727
$ws.Range("K40").Value = 2775.2727
$ws.Range("M40").Value = -2639.2727
# Row 46
$ws.Range("H46").Value = 4631.3335
$ws.Range("J46").Value = 4925
$ws.Range("L46").Value = 4925
$ws.Range("N46").Value = -5301
# Row 74
$ws.Range("H74").Value = 44288
$ws.Range("J74").Value = 56432
$ws.Range("L74").Value = 56432
$ws.Range("N74").Value = -58428
# Row 77
$ws.Range("H77").Value = 44288
$ws.Range("J77").Value = 56432
$ws.Range("L77").Value = 169296
$ws.Range("N77").Value = -179280
# Row 95
$ws.Range("H95").Value = 27194.5
$ws.Range("J95").Value = 27194.5
$ws.Range("L95").Value = 27194.5
$ws.Range("N95").Value = -32686.5
# Row 122
$ws.Range("H122").Value = 5387.6
$ws.Range("I122").Value = 4511.9414
$ws.Range("K122").Value = 13535.8242
$ws.Range("M122").Value = -11085.8242
# Row 136
$ws.Range("H136").Value = 63330.05
$ws.Range("I136").Value = 2569.7273
$ws.Range("J136").Value = 146875.5
$ws.Range("K136").Value = 7709.1819
$ws.Range("L136").Value = 440626.5
$ws.Range("M136").Value = -5159.1819
$ws.Range("N136").Value = -445726.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3665.6667
$ws.Range("I122").Value = 1671.8064
$ws.Range("K122").Value = 5015.4192
$ws.Range("M122").Value = -2565.4192
